$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Form Responses 1")

# Add the "Subject Number" names collected in column J for each survey
# response row (rows with no name were left blank by the author).
$names = [ordered]@{
    3  = " Alex Gustavo Chalco Maza"
    4  = "Usman Arshad"
    5  = "Destiny Berisha"
    6  = "Nadia Sultana"
    7  = "Fahmida Ferdousi"
    8  = "Hadia Perez"
    9  = "Jung Sang Cho"
    10 = "Kevin Call"
    11 = "Nicole Vazquez"
    12 = "James Castro"
    13 = "Jasmine Bachtarzi"
    14 = "Jean Carlos Huang"
    20 = "Taehyuk Kim"
    21 = "Syed Rizvi"
    22 = "Revital Schechter"
    23 = "Nicholas Carter"
    24 = "Aren Mineo"
    25 = "Ridmila Sudasinghe"
}

foreach ($row in $names.Keys) {
    $ws.Cells.Item($row, 10).Value = $names[$row]
}

$ws.Range("I7").Select() | Out-Null
